# B1 / B2 PowerPoint — Design change (Integral -> Office Theme) + table style refresh
#
# The deck's Design ("Integral", Red-Violet colours) is swapped for the
# default "Office Theme" colour scheme, and the table on slide 5 is
# re-styled to match (its a:tableStyleId changes from the custom
# "Table_0" GUID to the built-in medium-style-2-accent GUID that PowerPoint
# assigns once a new theme/style is applied).

$p = $ppt.ActivePresentation

# --- 1. Re-colour the presentation theme (Integral/Red Violet -> Office) ---
# All slides share one theme; grab it through any slide's ThemeColorScheme.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72

# --- 2. Point the slide-5 table at the new theme's table style ---
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{4A206071-D8CD-4B88-8B6D-C3E82DBDD06B}")
